# Add the 2024 season's timing/duration rows to Sheet1, continuing the
# existing Year / Day.count / TimePer table (cols A:C), then leave the
# selection on the next empty cell below the new data (matches the
# "Add files via upload" resave: table grows from row 100 to row 103).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Column A is text-formatted (style "3" from the <col> default), which would
# coerce a typed 2024 into a text value. Resetting the cell to the "Normal"
# style first (as Excel does when a fresh value overflows past the last
# styled row) keeps these new cells on the default/general style, so the
# year lands as a genuine number like the rest of the workbook expects.
$ws.Cells.Item(101, 1).Style = "Normal"
$ws.Cells.Item(101, 1).Value = 2024
$ws.Cells.Item(101, 2).Value = 197
$ws.Cells.Item(101, 3).Value = "Before"

$ws.Cells.Item(102, 1).Style = "Normal"
$ws.Cells.Item(102, 1).Value = 2024
$ws.Cells.Item(102, 2).Value = 38
$ws.Cells.Item(102, 3).Value = "During"

$ws.Cells.Item(103, 1).Style = "Normal"
$ws.Cells.Item(103, 1).Value = 2024
$ws.Cells.Item(103, 2).Value = 130
$ws.Cells.Item(103, 3).Value = "After"

# Land the selection just past the new data, same as the saved workbook.
$ws.Range("B104").Select()
